$wb = $excel.ActiveWorkbook

# Column-width helper constants: the report generator widens the "Status"
# style columns to fit "Handed back: in sync with en-US" and clamps any
# content-driven column at a 40-character cap. These ColumnWidth inputs are
# the closest settable values to that generator's output widths.
$statusColWidth = 29.166666666666668
$cappedColWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# 1) Overview sheet: the "zh-cn" / "de-de" status columns (E2/F2) move from
#    "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: handback completed - status flips, the "Latest Target File"
#    cell becomes a hyperlink to the source .md, the "Latest Handback File"
#    is populated with the generated .xlf, and the handback timestamp is set.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"

$zhHandoffAddress = ""
foreach ($hl in $wsZh.Hyperlinks) {
    $zhHandoffAddress = $hl.Address()
}
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhHandoffAddress, "", "", "cc54d45f-e499-49a1-8ed4-e1bdfd642a56.md")

$wsZh.Range("J2").Value = "cc54d45f-e499-49a1-8ed4-e1bdfd642a56.feba180e7bd4ad3af8137b31740a74d5d3a92386.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-23 00:55:47"

$wsZh.Columns.Item(3).ColumnWidth = $statusColWidth
$wsZh.Columns.Item(9).ColumnWidth = $cappedColWidth
$wsZh.Columns.Item(10).ColumnWidth = $cappedColWidth

# ---------------------------------------------------------------------------
# 3) de-de sheet: same handback shape, later timestamp (processed after zh-cn).
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"

$deHandoffAddress = ""
foreach ($hl in $wsDe.Hyperlinks) {
    $deHandoffAddress = $hl.Address()
}
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deHandoffAddress, "", "", "cc54d45f-e499-49a1-8ed4-e1bdfd642a56.md")

$wsDe.Range("J2").Value = "cc54d45f-e499-49a1-8ed4-e1bdfd642a56.feba180e7bd4ad3af8137b31740a74d5d3a92386.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-23 00:55:54"

$wsDe.Columns.Item(3).ColumnWidth = $statusColWidth
$wsDe.Columns.Item(9).ColumnWidth = $cappedColWidth
$wsDe.Columns.Item(10).ColumnWidth = $cappedColWidth
